# Generate Report for Handoff
# The localization status report is regenerated: the "In Translation" rows
# move to "Ready for handoff" and the handoff timestamps advance.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Status moves from "In Translation" to "Ready for handoff" everywhere it appears.
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Per-language "Latest Handoff Datetime" advances to the new handoff run time.
$zhcn.Range("H2").Value = "2016-08-26 22:38:30"
$dede.Range("H2").Value = "2016-08-26 22:38:35"

# Overview's "Latest HO Xliff Generate Date" advances to match the de-de handoff.
$overview.Range("G2").Value = "2016-08-26 22:38:35"

# The Status column is now wider because "Ready for handoff" is longer than
# "In Translation" -- mirror the resulting column widths.
$overview.Range("E1").ColumnWidth = 16.3
$overview.Range("F1").ColumnWidth = 16.3
$zhcn.Range("C1").ColumnWidth = 16.3
$dede.Range("C1").ColumnWidth = 16.3
